# MeanSDiffthdRanalysis.xlsx edit
# - Add two new columns (G "Threshold", H "noobject") with raw threshold data
#   and the "old" per-pair difference formulas (D column = G-H per matching
#   Observation/Group/Room/Duration row).
# - Rename the old "Threshold" header in E1 to "oldThreshold" (a second,
#   distinct shared string) while G1 keeps using the original "Threshold"
#   shared string.
# - Clear out a bunch of stray empty-but-styled cells (columns H:Q) that were
#   left over from earlier formatting and are no longer part of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Set G1 first so the original "Threshold" shared string stays referenced
# (and therefore alive) once E1 is renamed away from it.
$ws.Range("G1").Value = "Threshold"
$ws.Range("E1").Value = "oldThreshold"
$ws.Range("H1").Value = "noobject"
$ws.Range("H1").HorizontalAlignment = -4108

# --- New raw data in columns G (new Threshold) and H (noobject) -------
$ws.Range("G2").Value = 1.99
$ws.Range("H2").Value = 1.9

$ws.Range("G3").Value = 2.02
$ws.Range("H3").Value = 1.9

$ws.Range("G4").Value = 1.97
$ws.Range("H4").Value = 1.9

$ws.Range("G5").Value = 1.98
$ws.Range("H5").Value = 1.9

$ws.Range("G6").Value = 1.94
$ws.Range("H6").Value = 1.8699999999999999

$ws.Range("G7").Value = 1.98
$ws.Range("H7").Value = 1.8699999999999999

$ws.Range("G8").Value = 2.02
$ws.Range("H8").Value = 1.9799999999999998

$ws.Range("G9").Value = 2.0299999999999998
$ws.Range("H9").Value = 1.9799999999999998

$ws.Range("G10").Value = 1.94
$ws.Range("H10").Value = 1.89

$ws.Range("G11").Value = 1.94
$ws.Range("H11").Value = 1.89

$ws.Range("G12").Value = 1.97
$ws.Range("H12").Value = 1.94

$ws.Range("G13").Value = 1.98
$ws.Range("H13").Value = 1.94

# --- Difference formulas in column G for rows 14-25 --------------------
# Row 14+k mirrors row 2+k (same Observation/Group/Room/Duration), so the
# "old" difference is G(2+k)-H(2+k). Row 14 is a single formula; 15:25 is
# filled as one shared-formula block (matches a drag-fill in Excel).
$ws.Range("G14").Formula = "=G2-H2"
$ws.Range("G15:G25").Formula = "=G3-H3"

# --- Clear leftover stray formatted (empty) cells ----------------------
# (Issued as individual contiguous-range Clear() calls; multi-area Union
# ranges only clear their last area in this host, so keep these separate.)
$ws.Range("I13:Q13").Clear()
$ws.Range("P14:Q14").Clear()
$ws.Range("Q15").Clear()
$ws.Range("H16:Q16").Clear()
$ws.Range("H17:I17").Clear()
$ws.Range("P17:Q17").Clear()
$ws.Range("H18:I18").Clear()
$ws.Range("P18:Q18").Clear()
$ws.Range("H19:Q19").Clear()
$ws.Range("H20:J20").Clear()
$ws.Range("Q20").Clear()
$ws.Range("H21").Clear()
$ws.Range("H22").Clear()
